{"js": "// Update the date line and the twenty-five \"three digit x one digit\"\n// multiplication prompts in the table to the new day's values.\nconst replacements = [\n  [\"2024-03-19 Tuesday\", \"2024-03-20 Wednesday\"],\n  [\"517\u00d72=\", \"366\u00d78=\"],\n  [\"708\u00d79=\", \"497\u00d73=\"],\n  [\"620\u00d79=\", \"821\u00d74=\"],\n  [\"827\u00d79=\", \"185\u00d74=\"],\n  [\"772\u00d72=\", \"233\u00d76=\"],\n  [\"218\u00d79=\", \"169\u00d77=\"],\n  [\"421\u00d74=\", \"781\u00d79=\"],\n  [\"358\u00d77=\", \"654\u00d79=\"],\n  [\"187\u00d79=\", \"762\u00d76=\"],\n  [\"956\u00d75=\", \"398\u00d73=\"],\n  [\"421\u00d73=\", \"226\u00d77=\"],\n  [\"573\u00d77=\", \"335\u00d73=\"],\n  [\"745\u00d79=\", \"168\u00d76=\"],\n  [\"694\u00d77=\", \"928\u00d74=\"],\n  [\"318\u00d75=\", \"790\u00d78=\"],\n  [\"388\u00d75=\", \"987\u00d77=\"],\n  [\"307\u00d76=\", \"418\u00d72=\"],\n  [\"132\u00d77=\", \"432\u00d74=\"],\n  [\"367\u00d76=\", \"520\u00d72=\"],\n  [\"455\u00d79=\", \"150\u00d78=\"],\n  [\"303\u00d79=\", \"760\u00d77=\"],\n  [\"886\u00d79=\", \"135\u00d74=\"],\n  [\"593\u00d75=\", \"276\u00d73=\"],\n  [\"890\u00d74=\", \"857\u00d72=\"],\n  [\"688\u00d75=\", \"536\u00d73=\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2024-03-19 Tuesday\", \"2024-03-20 Wednesday\"),\n    @(\"517\u00d72=\", \"366\u00d78=\"),\n    @(\"708\u00d79=\", \"497\u00d73=\"),\n    @(\"620\u00d79=\", \"821\u00d74=\"),\n    @(\"827\u00d79=\", \"185\u00d74=\"),\n    @(\"772\u00d72=\", \"233\u00d76=\"),\n    @(\"218\u00d79=\", \"169\u00d77=\"),\n    @(\"421\u00d74=\", \"781\u00d79=\"),\n    @(\"358\u00d77=\", \"654\u00d79=\"),\n    @(\"187\u00d79=\", \"762\u00d76=\"),\n    @(\"956\u00d75=\", \"398\u00d73=\"),\n    @(\"421\u00d73=\", \"226\u00d77=\"),\n    @(\"573\u00d77=\", \"335\u00d73=\"),\n    @(\"745\u00d79=\", \"168\u00d76=\"),\n    @(\"694\u00d77=\", \"928\u00d74=\"),\n    @(\"318\u00d75=\", \"790\u00d78=\"),\n    @(\"388\u00d75=\", \"987\u00d77=\"),\n    @(\"307\u00d76=\", \"418\u00d72=\"),\n    @(\"132\u00d77=\", \"432\u00d74=\"),\n    @(\"367\u00d76=\", \"520\u00d72=\"),\n    @(\"455\u00d79=\", \"150\u00d78=\"),\n    @(\"303\u00d79=\", \"760\u00d77=\"),\n    @(\"886\u00d79=\", \"135\u00d74=\"),\n    @(\"593\u00d75=\", \"276\u00d73=\"),\n    @(\"890\u00d74=\", \"857\u00d72=\"),\n    @(\"688\u00d75=\", \"536\u00d73=\"),\n)\n\nforeach ($pair in $replacements) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $pair[0]\n    $find.Replacement.Text = $pair[1]\n    $find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n}\n"}
